$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Professional summary paragraph: neutralize "all Black and
#    Asian-American voters" -> "50M voters" (plain text, no formatting).
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# -----------------------------------------------------------------------
# 2) Siege Analytics bullet: same neutralization, but "50M" must stay
#    bold + colored like the "23%" / "64%" figures in the same bullet.
#    Replace only "all Black and Asian-American" (leave "voters," intact)
#    so the trailing text is untouched, then bold/color the new range.
# -----------------------------------------------------------------------
$bulletRange = $d.Content
$bulletRange.Find.Execute(
    "all Black and Asian-American", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$bulletRange.Text = "50M"
$bulletRange.Bold = 1
$bulletRange.Font.Color = 5258796   # wdColor BGR for RGB 2C3E50

# -----------------------------------------------------------------------
# 3) Key Projects "Impact" line: same neutralization (plain text).
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, improved electoral",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters nationwide, improved electoral", 2) | Out-Null

# -----------------------------------------------------------------------
# 4) Relocate the "Field Director - The Feldman Group" role block so it
#    appears immediately before "Software Engineer - Salsa Labs" (it
#    currently sits after "Programmer - Lake Research Partners", right
#    before "KEY PROJECTS").
# -----------------------------------------------------------------------
$fieldDirectorTexts = @(
    "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012",
    "Political Campaign Management",
    "• Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control",
    "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings",
    "• Created custom reports and data visualizations based on specific client requirements"
)

# 4a) Insert a fresh copy of the block right before the Salsa Labs heading.
$salsaFind = $d.Content
$salsaFind.Find.Execute(
    "Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$salsaHeadingPara = $salsaFind.Paragraphs(1)
$anchorPara = $salsaHeadingPara.Previous()
$anchorIndex = $anchorPara.Range.Paragraphs(1).Index

$insertPoint = $anchorPara.Range
$insertPoint.Collapse(0)   # wdCollapseEnd
foreach ($lineText in $fieldDirectorTexts) {
    $insertPoint.InsertParagraphAfter() | Out-Null
    $insertPoint.Collapse(0)
}

for ($i = 0; $i -lt $fieldDirectorTexts.Count; $i++) {
    $paraIndex = $anchorIndex + 1 + $i
    $d.Paragraphs($paraIndex).Range.Text = $fieldDirectorTexts[$i]
}
# Only the role heading carries the Heading 3 style; the rest inherit
# the Normal style from the anchor paragraph automatically.
$d.Paragraphs($anchorIndex + 1).Style = "Heading 3"

# 4b) Delete the original block (now the *second* occurrence of the
#     heading text, further down the document).
$oldFind = $d.Content
$oldFind.Find.Execute(
    "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oldFind.Collapse(0)   # wdCollapseEnd - skip past the copy we just inserted
$oldFind.Find.Execute(
    "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$oldHeadingPara = $oldFind.Paragraphs(1)
$oldStart = $oldHeadingPara.Range.Start
$oldLastPara = $oldHeadingPara
for ($i = 1; $i -lt $fieldDirectorTexts.Count; $i++) {
    $oldLastPara = $oldLastPara.Next()
}
$oldEnd = $oldLastPara.Range.End
$d.Range($oldStart, $oldEnd).Delete() | Out-Null
